# Updates cryptos list values (Price and Volume(1h) columns) per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.716.24"
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").Value = "'1.900.85"
$ws.Range("E3").Value = '  +0.28%  '
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("E5").Value = '  +0.04%  '
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").Value = "'0.5229"
$ws.Range("E7").Value = '  +6.15%  '
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").Value = "'0.07247"
$ws.Range("E9").Value = '  -1.06%  '
$ws.Range("E10").Value = '  +2.89%  '
$ws.Range("D11").Value = "'0.9026"
$ws.Range("E11").Value = '  -0.88%  '
$ws.Range("E12").Value = '  +0.30%  '
$ws.Range("D13").Value = "'1.919.75"
$ws.Range("E13").Value = '  +1.41%  '
$ws.Range("E14").Value = '  -0.22%  '
$ws.Range("D15").Value = "'92.32"
$ws.Range("E15").Value = '  +1.35%  '
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("D17").Value = "'0.000008699"
$ws.Range("E17").Value = '  -0.37%  '
$ws.Range("D18").Value = "'1.000"
$ws.Range("E18").Value = '  +0.10%  '
$ws.Range("D19").Value = "'27.752.49"
$ws.Range("E19").Value = '  -0.19%  '
$ws.Range("D20").Value = "'14.46"
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("E21").Value = '  +0.51%  '
$ws.Range("D22").Value = "'2.141.92"
$ws.Range("E22").Value = '  +1.32%  '
$ws.Range("E23").Value = '  +0.93%  '
$ws.Range("D24").Value = "'6.611"
$ws.Range("E24").Value = '  -0.60%  '
$ws.Range("E25").Value = '  -0.49%  '
$ws.Range("D26").Value = "'1.867"
$ws.Range("E26").Value = '  +0.99%  '
$ws.Range("D27").Value = "'18.31"
$ws.Range("E27").Value = '  -0.48%  '
$ws.Range("D28").Value = "'2.161"
$ws.Range("E28").Value = '  -0.71%  '
$ws.Range("D29").Value = "'114.53"
$ws.Range("E29").Value = '  -0.64%  '
$ws.Range("D30").Value = "'4.846"
$ws.Range("E30").Value = '  -0.69%  '
$ws.Range("D31").Value = "'0.09092"
$ws.Range("E31").Value = '  +1.78%  '
$ws.Range("D32").Value = "'3.189"
$ws.Range("E32").Value = '  -1.19%  '
$ws.Range("D33").Value = "'4.837"
$ws.Range("E33").Value = '  +4.09%  '
$ws.Range("D34").Value = "'1.225"
$ws.Range("E34").Value = '  -0.33%  '
$ws.Range("D35").Value = "'0.7793"
$ws.Range("E35").Value = '  +1.62%  '
$ws.Range("D36").Value = "'0.02093"
$ws.Range("E36").Value = '  +2.35%  '
$ws.Range("D37").Value = "'2.577"
$ws.Range("E37").Value = '  +0.75%  '
$ws.Range("E38").Value = '  +2.66%  '
$ws.Range("D39").Value = "'1.094"
$ws.Range("E39").Value = '  -0.53%  '
$ws.Range("D40").Value = "'0.5548"
$ws.Range("E40").Value = '  +0.96%  '
$ws.Range("D41").Value = "'0.05290"
$ws.Range("E41").Value = '  +0.15%  '
$ws.Range("D42").Value = "'6.730"
$ws.Range("E42").Value = '  -2.44%  '
$ws.Range("D43").Value = "'116.14"
$ws.Range("E43").Value = '  +3.06%  '
$ws.Range("D44").Value = "'8.511"
$ws.Range("E44").Value = '  -0.62%  '
$ws.Range("D45").Value = "'0.1517"
$ws.Range("E45").Value = '  -0.17%  '
$ws.Range("D46").Value = "'0.4814"
$ws.Range("E46").Value = '  +0.44%  '
$ws.Range("D47").Value = "'10.48"
$ws.Range("E47").Value = '  -0.73%  '
$ws.Range("E48").Value = '  +0.12%  '
$ws.Range("D49").Value = "'1.613"
$ws.Range("E49").Value = '  -1.18%  '
$ws.Range("D50").Value = "'66.76"
$ws.Range("E50").Value = '  -0.96%  '
$ws.Range("D51").Value = "'0.06003"
$ws.Range("E51").Value = '  -0.90%  '
